# Update Col4a1-Itgav LR-pair sheet with recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Ligand-side values (columns G,H,I,J) keyed by Sending cluster (column A)
$ghij = @{
    "ECs" = @(237.699173, 475.398346, 0.6281175046907784, 0.5910692918834948)
    "FAPs" = @(46.00982733333333, 138.029482, 0.1215804732137916, 0.1716139504296584)
    "Inflammatory-Mac" = @(0.424496, 1.273488, 0.001121726108282323, 0.001583345118289764)
    "MuSCs" = @(93.29178250000001, 186.583565, 0.2465225305266621, 0.2319819085816677)
    "Neutrophils" = @(0.7462629999999999, 2.238789, 0.001971991940430749, 0.002783517107370325)
    "Resolving-Mac" = @(0.259518, 0.778554, 0.0006857735200548696, 0.0009679868795190599)
}

# New Receptor-side values (columns M,N,O,P) keyed by Target cluster (column D)
$mnop = @{
    "ECs" = @(20.574342, 41.148684, 0.07442291871210138, 0.05295769307665528)
    "FAPs" = @(58.255493, 174.766479, 0.2107257583291067, 0.2249216416002446)
    "Inflammatory-Mac" = @(59.53576899999999, 178.607307, 0.2153568603433074, 0.2298647253300724)
    "MuSCs" = @(31.770234, 63.540468, 0.1149214658940947, 0.08177555817559162)
    "Neutrophils" = @(31.46548433333334, 94.39645300000001, 0.1138191045320704, 0.121486713536183)
    "Resolving-Mac" = @(74.85037233333333, 224.551117, 0.2707538921893195, 0.2889936682812531)
}

# New Edge values (columns Q,R,S,T) keyed by row number (2..37)
$qrst = @{
    2 = @(4890.504078419166, 19562.01631367666, 0.04674633799324976, 0.03130166614660209)
    3 = @(13847.28250880729, 83083.69505284373, 0.1323605374957505, 0.1329442754299298)
    4 = @(14151.60305521904, 84909.61833131422, 0.1352694137368787, 0.1358659804298399)
    5 = @(7551.758347816482, 30207.03339126593, 0.07218418439280513, 0.04833502126422447)
    6 = @(7479.31960407779, 44875.91762446674, 0.07149177192482292, 0.07180706574308464)
    7 = @(17791.87160237541, 106751.2296142525, 0.1700652591472714, 0.1708152828698139)
    8 = @(946.621922916948, 5679.731537501689, 0.00904837367496883, 0.009088278914526182)
    9 = @(2680.325174148209, 24122.92656733388, 0.02562013741598787, 0.03859969145214177)
    10 = @(2739.230451847219, 24653.07406662497, 0.02618318899037574, 0.03944799357832208)
    11 = @(1461.742980679596, 8770.457884077578, 0.01397220620582664, 0.01403382658710363)
    12 = @(1447.721501136372, 13029.49351022735, 0.01383818058977913, 0.02084881483466061)
    13 = @(3443.852706892377, 30994.67436203139, 0.03291838633685337, 0.04959534506290411)
    14 = @(8.733725881632001, 52.402355289792, 0.00008348213097393717, 0.00008385030480880975)
    15 = @(24.729223756528, 222.563013808752, 0.0002363765848053503, 0.0003561285832254671)
    16 = @(25.272695797424, 227.454262176816, 0.000241571412844798, 0.0003639551907183874)
    17 = @(13.486337252064, 80.918023512384, 0.0001289104086954825, 0.0001294789308327435)
    18 = @(13.35697223756267, 120.212750138064, 0.0001276738611749383, 0.0001923553948145822)
    19 = @(31.77368365401066, 285.9631528860959, 0.0003037117097878171, 0.0004575767138897734)
    20 = @(1919.417038944615, 7677.668155778461, 0.0183469262500873, 0.01228522671400466)
    21 = @(5434.758782386273, 32608.55269431764, 0.05194864719044123, 0.05217775169974657)
    22 = @(5554.198012518243, 33325.18807510946, 0.05309031817810911, 0.053324457697671)
    23 = @(2963.901760302106, 11855.60704120842, 0.02833073058404571, 0.01897045006090494)
    24 = @(2935.471120682491, 17612.82672409495, 0.02805897367152467, 0.02818271967343805)
    25 = @(6982.924655765351, 41897.5479345921, 0.0667469346524541, 0.06704130273590246)
    26 = @(15.353870183946, 92.12322110367599, 0.0001467613958835967, 0.000147408644645737)
    27 = @(43.47391897265899, 391.265270753931, 0.0004155494970661563, 0.0006260732372120979)
    28 = @(44.42934158124699, 399.8640742312229, 0.0004246819929134726, 0.0006398323953372373)
    29 = @(23.708950135542, 142.253700813252, 0.0002266242045256419, 0.0002276236651465166)
    30 = @(23.48152673504633, 211.333740615417, 0.0002244503568042878, 0.0003381603454461633)
    31 = @(55.85806340859032, 502.7225706773129, 0.000533924493237594, 0.000804418819582573)
    32 = @(5.339412087156, 32.036472522936, 0.00005103726693795519, 0.00005126235206779967)
    33 = @(15.118349032374, 136.065141291366, 0.0001445101450555833, 0.0002177211979889252)
    34 = @(15.450603699342, 139.055433294078, 0.0001476860321855949, 0.0002225060381837625)
    35 = @(8.244947587212, 49.469685523272, 0.00007881009819605894, 0.00007915766737932028)
    36 = @(8.165859563218, 73.492736068962, 0.00007805412796445107, 0.0001175975447389157)
    37 = @(19.425018927202, 174.825170344818, 0.0001856758497152263, 0.0002797420791603365)
}

for ($r = 2; $r -le 37; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target  = $ws.Cells.Item($r, 4).Value2

    $g = $ghij[$sending]
    $ws.Cells.Item($r, 7).Value  = $g[0]
    $ws.Cells.Item($r, 8).Value  = $g[1]
    $ws.Cells.Item($r, 9).Value  = $g[2]
    $ws.Cells.Item($r, 10).Value = $g[3]

    $m = $mnop[$target]
    $ws.Cells.Item($r, 13).Value = $m[0]
    $ws.Cells.Item($r, 14).Value = $m[1]
    $ws.Cells.Item($r, 15).Value = $m[2]
    $ws.Cells.Item($r, 16).Value = $m[3]

    $q = $qrst[$r]
    $ws.Cells.Item($r, 17).Value = $q[0]
    $ws.Cells.Item($r, 18).Value = $q[1]
    $ws.Cells.Item($r, 19).Value = $q[2]
    $ws.Cells.Item($r, 20).Value = $q[3]
}
